$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:C2").Copy()
$ws.Range("A78:C78").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A79:C79").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A80").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B80").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C80").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A81").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B81").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C81").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B82").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C82").PasteSpecial(-4122)

$ws.Range("A67:C67").Copy()
$ws.Range("A83:C83").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A84:C84").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A85").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B85").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C85").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A86").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B86").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C86").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A87").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B87").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C87").PasteSpecial(-4122)

$ws.Range("A2:C2").Copy()
$ws.Range("A88:C88").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A89:C89").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A90").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B90").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C90").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B91").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C91").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A92").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B92").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C92").PasteSpecial(-4122)

$ws.Range("A67:C67").Copy()
$ws.Range("A93:C93").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A94:C94").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A95").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B95").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C95").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A96").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B96").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C96").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A97").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B97").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C97").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A98:C98").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A99").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B99").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C99").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A100").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B100").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C100").PasteSpecial(-4122)

$ws.Range("A2:C2").Copy()
$ws.Range("A101:C101").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A102:C102").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A103").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B103").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C103").PasteSpecial(-4122)

$ws.Range("A2:C2").Copy()
$ws.Range("A104:C104").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A105").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B105").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C105").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A106").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B106").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C106").PasteSpecial(-4122)

$ws.Range("A67:C67").Copy()
$ws.Range("A107:C107").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A108:C108").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A109").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B109").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C109").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A110").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B110").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C110").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A111").PasteSpecial(-4122)
$ws.Range("B68").Copy()
$ws.Range("B111").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C111").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A112:C112").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("A78:C78").Merge()
$ws.Range("A88:C88").Merge()
$ws.Range("A101:C101").Merge()
$ws.Range("A104:C104").Merge()

$ws.Range("A78").Value = "REST Services & Consuming the REST Services using Postman"
$ws.Range("A79").Value = "REST Introduction"
$ws.Range("C79").Value = 43459
$ws.Range("B80").Value = "Introduction"
$ws.Range("B81").Value = "REST Vs SOAP"
$ws.Range("B82").Value = "JSON Structure"
$ws.Range("B83").Value = "Installaing POSTMAN"
$ws.Range("A84").Value = "Consuming REST Services"
$ws.Range("C84").Value = 43460
$ws.Range("B85").Value = "POST, GET, PUT, DELETE methods"
$ws.Range("B86").Value = "HTTP Response Codes"
$ws.Range("B87").Value = "Demos"
$ws.Range("A88").Value = "Java Basics, Collections and Java 8"
$ws.Range("A89").Value = "Java Basics"
$ws.Range("C89").Value = 43461
$ws.Range("B90").Value = "OOPS concepts"
$ws.Range("A91").Value = "Collections"
$ws.Range("C91").Value = 43462
$ws.Range("B92").Value = "Implements and Extends concepts"
$ws.Range("B93").Value = "Threading and Synchronization concepts"
$ws.Range("B94").Value = "Collections"
$ws.Range("A95").Value = "Maps"
$ws.Range("C95").Value = 43463
$ws.Range("B96").Value = "Collections Demo"
$ws.Range("B97").Value = "Maps & Demo"
$ws.Range("A98").Value = "Java 8"
$ws.Range("C98").Value = 43467
$ws.Range("B99").Value = "Object Oriented Vs Functional"
$ws.Range("B100").Value = "Java 8 Concepts"
$ws.Range("A101").Value = "UNIX"
$ws.Range("A102").Value = "Unix Commands"
$ws.Range("C102").Value = 43468
$ws.Range("B103").Value = "Commands"
$ws.Range("A104").Value = "Jenkins"
$ws.Range("A105").Value = "Jenkins Basics"
$ws.Range("C105").Value = 43469
$ws.Range("B106").Value = "Jenkins Introduction"
$ws.Range("B107").Value = "Jenkins installation"
$ws.Range("B108").Value = "Demo"
$ws.Range("A109").Value = "Jenkins Plugins"
$ws.Range("C109").Value = 43472
$ws.Range("B110").Value = "Configuring environments"
$ws.Range("B111").Value = "Plugins required"
$ws.Range("B112").Value = "Demo"
